$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray empty paragraph right after "Wilt advies voor website"
# ---------------------------------------------------------------------------
$rng = $d.Content
$ok = $rng.Find.Execute("Wilt advies voor website")
if (-not $ok) { throw "Could not find 'Wilt advies voor website'" }
$rng.Collapse(0)
$rng.MoveEnd(1, 1)
$rng.Delete()

# ---------------------------------------------------------------------------
# 2) Append a new run to the paragraph that ends with
#    "...alles duidelijk in beeld krijgen. " (append after the trailing space run)
# ---------------------------------------------------------------------------
$rng = $d.Content
$ok = $rng.Find.Execute("alles duidelijk in beeld krijgen.")
if (-not $ok) { throw "Could not find 'alles duidelijk in beeld krijgen.'" }
$rng.Expand(4)
$rng.Collapse(0)
$rng.MoveEnd(-1, 1)
$rng.InsertAfter(" Alle dingen die u in de site wilt hebben zal toegevoegd worden en er voor zorgen dat alles perfect werkt")

# ---------------------------------------------------------------------------
# 3) After "... en betrokken bent bij het maken van de site." add a new
#    paragraph made of two runs.
# ---------------------------------------------------------------------------
$rng = $d.Content
$ok = $rng.Find.Execute("en betrokken bent bij het maken van de site.")
if (-not $ok) { throw "Could not find 'en betrokken bent bij het maken van de site.'" }
$insertAt = $rng.End
$ip = $d.Range($insertAt, $insertAt)
$ip.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:t xml:space='preserve'>Wie bieden u snelle en goede service. Wij staan open voor verandering en reageren snel op reacties en issues. </w:t></w:r><w:r><w:t xml:space='preserve'> Filters zodat de klant alles kan bepalen en kan zien of het bezet is. In de filters worden Maximaal bedrag en minimaal bedrag, hoeveel personen en dergelijke benodigheden.</w:t></w:r></w:p>")

# ---------------------------------------------------------------------------
# 4) After "...Ook daar zult u dingen kunnen aanpassen en bent u erbij betrokken."
#    add a new paragraph (single run).
# ---------------------------------------------------------------------------
$rng = $d.Content
$ok = $rng.Find.Execute("Ook daar zult u dingen kunnen aanpassen en bent u erbij betrokken.")
if (-not $ok) { throw "Could not find 'Ook daar zult u dingen kunnen aanpassen en bent u erbij betrokken.'" }
$insertAt = $rng.End
$ip = $d.Range($insertAt, $insertAt)
$ip.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:t xml:space='preserve'>Wij hebben al eerdere ervaring met het maken van website met de vraag van de klant over zaal verhuur en kantoor verhuur. </w:t></w:r></w:p>")

# ---------------------------------------------------------------------------
# 5) Replace the paragraph "Ook een plattengrond van het gebouw toevoegen."
#    (two runs) with a single run with new text.
# ---------------------------------------------------------------------------
$rng = $d.Content
$ok = $rng.Find.Execute("Ook een plattengrond")
if (-not $ok) { throw "Could not find 'Ook een plattengrond'" }
$rng.Expand(4)
$rng.Text = "Wij kunnen ook een plattegrond toevoegen van de trechter om de consument te laten zien wat het is en hoe groot het is. We kunnen ook een live chat toevoegen tussen ons zodat u 24/7 contact met ons kunt opzoeken. "

# ---------------------------------------------------------------------------
# 6) Add a new paragraph after that one: "Wat ons uniek maakt..."
# ---------------------------------------------------------------------------
$rng = $d.Content
$ok = $rng.Find.Execute("live chat toevoegen tussen ons zodat u 24/7 contact met ons kunt opzoeken.")
if (-not $ok) { throw "Could not find the live-chat sentence" }
$insertAt = $rng.End
$ip = $d.Range($insertAt, $insertAt)
$ip.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:t xml:space='preserve'>Wat ons uniek maakt vergeleken met andere is dat we sites kunnen maken met bewegende achtergronden om echt het unieke beeld te creëren.  We doen er ook niet lang over!</w:t></w:r></w:p>")

# ---------------------------------------------------------------------------
# 7) Two trailing empty paragraphs at the very end of the document body.
# ---------------------------------------------------------------------------
$endPoint = $d.Range($d.Content.End, $d.Content.End)
$endPoint.InsertAfter("`r`r")

Write-Output "done"
